# Applies the commit: reworks the ODI Batting MATCH_CARD_LINK column into a
# MATCH_CODE column holding just the numeric match code, then inserts a new
# "Player Info" sheet ahead of "ODI Batting".

$wb = $excel.ActiveWorkbook
$batting = $wb.Worksheets.Item("ODI Batting")

# --- 1. Rework the ODI Batting "MATCH_CARD_LINK" column (do this first,
#        while $batting is still guaranteed to resolve to the right sheet -
#        inserting a sheet ahead of it later can otherwise leave a stale
#        reference behind). ---
$batting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @{
    2  = "4075"
    3  = "4285"
    4  = "4286"
    5  = "4291"
    6  = "4293"
    7  = "4296"
    8  = "4338"
    9  = "4344"
    10 = "4385"
    11 = "4397"
    12 = "4413"
    13 = "4414"
    14 = "4417"
    15 = "4443"
    16 = "4445"
    17 = "4447"
}

foreach ($row in $matchCodes.Keys) {
    $cell = $batting.Cells.Item($row, 4)
    # The workbook stores every value (even numeric-looking ones) as text,
    # so force text by prefixing with an apostrophe, then strip the
    # resulting quote-prefix style back off so the exported cell carries no
    # style id (matching the rest of the column).
    $cell.Value = "'" + $matchCodes[$row]
    $cell.Style = "Normal"
}

# --- 2. Insert the new "Player Info" worksheet ahead of "ODI Batting" ---
$info = $wb.Worksheets.Add($batting)
$info.Name = "Player Info"

# Header row
$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

$headerRange = $info.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# Data row - ID stays textual, same trick as above.
$info.Range("A2").Value = "'4702"
$info.Range("A2").Style = "Normal"
$info.Range("B2").Value = "Sunil Walford Ambris"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Does Not Bowl | Unknown"

$info.Range("A1").Select() | Out-Null

Write-Output "Applied Player Info sheet + MATCH_CODE rework"
